# Fill in newly-reported Intermediate1 (BY:CC), National4 (BJ:BN) and
# National3 (BO:BS) candidate/entry/pass counts for years 2000-2019
# across the three data sheets. Dash ("-") marks suppressed/NA figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("total_candidates")
$ws.Range("BY17").Value = 6695
$ws.Range("BZ17").Value = 3061
$ws.Range("CA17").Value = 3634
$ws.Range("CB17").Value = "-"
$ws.Range("CC17").Value = "-"
$ws.Range("BY18").Value = 10841
$ws.Range("BZ18").Value = 5096
$ws.Range("CA18").Value = 5745
$ws.Range("CB18").Value = "-"
$ws.Range("CC18").Value = "-"
$ws.Range("BY19").Value = 13925
$ws.Range("BZ19").Value = 6633
$ws.Range("CA19").Value = 7292
$ws.Range("CB19").Value = "-"
$ws.Range("CC19").Value = "-"
$ws.Range("BY20").Value = 17693
$ws.Range("BZ20").Value = 8665
$ws.Range("CA20").Value = 9028
$ws.Range("CB20").Value = "-"
$ws.Range("CC20").Value = "-"
$ws.Range("BY21").Value = 22720
$ws.Range("BZ21").Value = 11079
$ws.Range("CA21").Value = 11641
$ws.Range("CB21").Value = "-"
$ws.Range("CC21").Value = "-"
$ws.Range("BY22").Value = 26906
$ws.Range("BZ22").Value = 13057
$ws.Range("CA22").Value = 13849
$ws.Range("CB22").Value = "-"
$ws.Range("CC22").Value = "-"
$ws.Range("BY23").Value = 32352
$ws.Range("BZ23").Value = 15915
$ws.Range("CA23").Value = 16437
$ws.Range("CB23").Value = "-"
$ws.Range("CC23").Value = "-"
$ws.Range("BY24").Value = 37490
$ws.Range("BZ24").Value = 18494
$ws.Range("CA24").Value = 18996
$ws.Range("CB24").Value = "-"
$ws.Range("CC24").Value = "-"
$ws.Range("BY25").Value = 40693
$ws.Range("BZ25").Value = 20319
$ws.Range("CA25").Value = 20374
$ws.Range("CB25").Value = "-"
$ws.Range("CC25").Value = "-"
$ws.Range("BY26").Value = 43061
$ws.Range("BZ26").Value = 21534
$ws.Range("CA26").Value = 21527
$ws.Range("CB26").Value = "-"
$ws.Range("CC26").Value = "-"
$ws.Range("BY27").Value = 44837
$ws.Range("BZ27").Value = 22422
$ws.Range("CA27").Value = 22415
$ws.Range("CB27").Value = "-"
$ws.Range("CC27").Value = "-"
$ws.Range("BY28").Value = 45645
$ws.Range("BZ28").Value = 23132
$ws.Range("CA28").Value = 22513
$ws.Range("CB28").Value = "-"
$ws.Range("CC28").Value = "-"
$ws.Range("BY29").Value = 45882
$ws.Range("BZ29").Value = 23666
$ws.Range("CA29").Value = 22216
$ws.Range("CB29").Value = "-"
$ws.Range("CC29").Value = "-"
$ws.Range("CB30").Value = "-"
$ws.Range("CC30").Value = "-"
$ws.Range("BJ31").Value = 40641
$ws.Range("BK31").Value = 21537
$ws.Range("BL31").Value = 19104
$ws.Range("BM31").Value = "-"
$ws.Range("BN31").Value = "-"
$ws.Range("BO31").Value = 16650
$ws.Range("BP31").Value = 10135
$ws.Range("BQ31").Value = 6515
$ws.Range("BR31").Value = "-"
$ws.Range("BS31").Value = "-"
$ws.Range("BY31").Value = 12241
$ws.Range("BZ31").Value = 6775
$ws.Range("CA31").Value = 5466
$ws.Range("CB31").Value = "-"
$ws.Range("CC31").Value = "-"
$ws.Range("BJ32").Value = 53111
$ws.Range("BK32").Value = 28377
$ws.Range("BL32").Value = 24734
$ws.Range("BM32").Value = "-"
$ws.Range("BN32").Value = "-"
$ws.Range("BO32").Value = 10895
$ws.Range("BP32").Value = 6389
$ws.Range("BQ32").Value = 4506
$ws.Range("BR32").Value = "-"
$ws.Range("BS32").Value = "-"
$ws.Range("BY32").Value = 1394
$ws.Range("BZ32").Value = 708
$ws.Range("CA32").Value = 686
$ws.Range("CB32").Value = "-"
$ws.Range("CC32").Value = "-"
$ws.Range("BJ33").Value = 53219
$ws.Range("BK33").Value = 28499
$ws.Range("BL33").Value = 24720
$ws.Range("BM33").Value = "-"
$ws.Range("BN33").Value = "-"
$ws.Range("BO33").Value = 11223
$ws.Range("BP33").Value = 6529
$ws.Range("BQ33").Value = 4694
$ws.Range("BR33").Value = "-"
$ws.Range("BS33").Value = "-"
$ws.Range("BY33").Value = "-"
$ws.Range("BZ33").Value = "-"
$ws.Range("CA33").Value = "-"
$ws.Range("CB33").Value = "-"
$ws.Range("CC33").Value = "-"
$ws.Range("BJ34").Value = 51970
$ws.Range("BK34").Value = 27567
$ws.Range("BL34").Value = 24403
$ws.Range("BM34").Value = "-"
$ws.Range("BN34").Value = "-"
$ws.Range("BO34").Value = 10465
$ws.Range("BP34").Value = 6081
$ws.Range("BQ34").Value = 4384
$ws.Range("BR34").Value = "-"
$ws.Range("BS34").Value = "-"
$ws.Range("BY34").Value = "-"
$ws.Range("BZ34").Value = "-"
$ws.Range("CA34").Value = "-"
$ws.Range("CB34").Value = "-"
$ws.Range("CC34").Value = "-"
$ws.Range("BJ35").Value = 46774
$ws.Range("BK35").Value = 24763
$ws.Range("BL35").Value = 22011
$ws.Range("BM35").Value = "-"
$ws.Range("BN35").Value = "-"
$ws.Range("BO35").Value = 10795
$ws.Range("BP35").Value = 6354
$ws.Range("BQ35").Value = 4441
$ws.Range("BR35").Value = "-"
$ws.Range("BS35").Value = "-"
$ws.Range("BY35").Value = "-"
$ws.Range("BZ35").Value = "-"
$ws.Range("CA35").Value = "-"
$ws.Range("CB35").Value = "-"
$ws.Range("CC35").Value = "-"
$ws.Range("BJ36").Value = 46544
$ws.Range("BK36").Value = 25029
$ws.Range("BL36").Value = 21512
$ws.Range("BM36").Value = 3
$ws.Range("BN36").Value = "-"
$ws.Range("BO36").Value = 11380
$ws.Range("BP36").Value = 6734
$ws.Range("BQ36").Value = 4645
$ws.Range("BR36").Value = 1
$ws.Range("BS36").Value = "-"
$ws.Range("BY36").Value = "-"
$ws.Range("BZ36").Value = "-"
$ws.Range("CA36").Value = "-"
$ws.Range("CB36").Value = "-"
$ws.Range("CC36").Value = "-"

$ws = $wb.Worksheets.Item("total_subject_entries")
$ws.Range("BY17").Value = 8096
$ws.Range("BZ17").Value = 3706
$ws.Range("CA17").Value = 4390
$ws.Range("CB17").Value = "-"
$ws.Range("CC17").Value = "-"
$ws.Range("BY18").Value = 15981
$ws.Range("BZ18").Value = 7534
$ws.Range("CA18").Value = 8447
$ws.Range("CB18").Value = "-"
$ws.Range("CC18").Value = "-"
$ws.Range("BY19").Value = 20352
$ws.Range("BZ19").Value = 9826
$ws.Range("CA19").Value = 10526
$ws.Range("CB19").Value = "-"
$ws.Range("CC19").Value = "-"
$ws.Range("BY20").Value = 24613
$ws.Range("BZ20").Value = 12157
$ws.Range("CA20").Value = 12456
$ws.Range("CB20").Value = "-"
$ws.Range("CC20").Value = "-"
$ws.Range("BY21").Value = 31231
$ws.Range("BZ21").Value = 15168
$ws.Range("CA21").Value = 16063
$ws.Range("CB21").Value = "-"
$ws.Range("CC21").Value = "-"
$ws.Range("BY22").Value = 36653
$ws.Range("BZ22").Value = 17678
$ws.Range("CA22").Value = 18975
$ws.Range("CB22").Value = "-"
$ws.Range("CC22").Value = "-"
$ws.Range("BY23").Value = 45174
$ws.Range("BZ23").Value = 22323
$ws.Range("CA23").Value = 22851
$ws.Range("CB23").Value = "-"
$ws.Range("CC23").Value = "-"
$ws.Range("BY24").Value = 53840
$ws.Range("BZ24").Value = 26455
$ws.Range("CA24").Value = 27385
$ws.Range("CB24").Value = "-"
$ws.Range("CC24").Value = "-"
$ws.Range("BY25").Value = 60267
$ws.Range("BZ25").Value = 29769
$ws.Range("CA25").Value = 30498
$ws.Range("CB25").Value = "-"
$ws.Range("CC25").Value = "-"
$ws.Range("BY26").Value = 65735
$ws.Range("BZ26").Value = 32752
$ws.Range("CA26").Value = 32983
$ws.Range("CB26").Value = "-"
$ws.Range("CC26").Value = "-"
$ws.Range("BY27").Value = 69834
$ws.Range("BZ27").Value = 34870
$ws.Range("CA27").Value = 34964
$ws.Range("CB27").Value = "-"
$ws.Range("CC27").Value = "-"
$ws.Range("BY28").Value = 72324
$ws.Range("BZ28").Value = 36585
$ws.Range("CA28").Value = 35739
$ws.Range("CB28").Value = "-"
$ws.Range("CC28").Value = "-"
$ws.Range("BY29").Value = 74383
$ws.Range("BZ29").Value = 38157
$ws.Range("CA29").Value = 36226
$ws.Range("CB29").Value = "-"
$ws.Range("CC29").Value = "-"
$ws.Range("CB30").Value = "-"
$ws.Range("CC30").Value = "-"
$ws.Range("BJ31").Value = 131711
$ws.Range("BK31").Value = 72910
$ws.Range("BL31").Value = 58801
$ws.Range("BM31").Value = "-"
$ws.Range("BN31").Value = "-"
$ws.Range("BO31").Value = 10240
$ws.Range("BP31").Value = 6035
$ws.Range("BQ31").Value = 4205
$ws.Range("BR31").Value = "-"
$ws.Range("BS31").Value = "-"
$ws.Range("BY31").Value = 18679
$ws.Range("BZ31").Value = 10584
$ws.Range("CA31").Value = 8095
$ws.Range("CB31").Value = "-"
$ws.Range("CC31").Value = "-"
$ws.Range("BJ32").Value = 140151
$ws.Range("BK32").Value = 78441
$ws.Range("BL32").Value = 61710
$ws.Range("BM32").Value = "-"
$ws.Range("BN32").Value = "-"
$ws.Range("BO32").Value = 17525
$ws.Range("BP32").Value = 10455
$ws.Range("BQ32").Value = 7070
$ws.Range("BR32").Value = "-"
$ws.Range("BS32").Value = "-"
$ws.Range("BY32").Value = 1553
$ws.Range("BZ32").Value = 792
$ws.Range("CA32").Value = 761
$ws.Range("CB32").Value = "-"
$ws.Range("CC32").Value = "-"
$ws.Range("BJ33").Value = 131747
$ws.Range("BK33").Value = 73643
$ws.Range("BL33").Value = 58104
$ws.Range("BM33").Value = "-"
$ws.Range("BN33").Value = "-"
$ws.Range("BO33").Value = 18596
$ws.Range("BP33").Value = 11292
$ws.Range("BQ33").Value = 7304
$ws.Range("BR33").Value = "-"
$ws.Range("BS33").Value = "-"
$ws.Range("BY33").Value = "-"
$ws.Range("BZ33").Value = "-"
$ws.Range("CA33").Value = "-"
$ws.Range("CB33").Value = "-"
$ws.Range("CC33").Value = "-"
$ws.Range("BJ34").Value = 125649
$ws.Range("BK34").Value = 70188
$ws.Range("BL34").Value = 55461
$ws.Range("BM34").Value = "-"
$ws.Range("BN34").Value = "-"
$ws.Range("BO34").Value = 17282
$ws.Range("BP34").Value = 10360
$ws.Range("BQ34").Value = 6922
$ws.Range("BR34").Value = "-"
$ws.Range("BS34").Value = "-"
$ws.Range("BY34").Value = "-"
$ws.Range("BZ34").Value = "-"
$ws.Range("CA34").Value = "-"
$ws.Range("CB34").Value = "-"
$ws.Range("CC34").Value = "-"
$ws.Range("BJ35").Value = 109644
$ws.Range("BK35").Value = 61471
$ws.Range("BL35").Value = 48173
$ws.Range("BM35").Value = "-"
$ws.Range("BN35").Value = "-"
$ws.Range("BO35").Value = 17905
$ws.Range("BP35").Value = 10845
$ws.Range("BQ35").Value = 7060
$ws.Range("BR35").Value = "-"
$ws.Range("BS35").Value = "-"
$ws.Range("BY35").Value = "-"
$ws.Range("BZ35").Value = "-"
$ws.Range("CA35").Value = "-"
$ws.Range("CB35").Value = "-"
$ws.Range("CC35").Value = "-"
$ws.Range("BJ36").Value = 107999
$ws.Range("BK36").Value = 61212
$ws.Range("BL36").Value = 46783
$ws.Range("BM36").Value = 4
$ws.Range("BN36").Value = "-"
$ws.Range("BO36").Value = 19285
$ws.Range("BP36").Value = 11804
$ws.Range("BQ36").Value = 7480
$ws.Range("BR36").Value = 1
$ws.Range("BS36").Value = "-"
$ws.Range("BY36").Value = "-"
$ws.Range("BZ36").Value = "-"
$ws.Range("CA36").Value = "-"
$ws.Range("CB36").Value = "-"
$ws.Range("CC36").Value = "-"

$ws = $wb.Worksheets.Item("total_passes")
$ws.Range("BY17").Value = 5804
$ws.Range("BZ17").Value = 2535
$ws.Range("CA17").Value = 3269
$ws.Range("CB17").Value = "-"
$ws.Range("CC17").Value = "-"
$ws.Range("BY18").Value = 10444
$ws.Range("BZ18").Value = 4706
$ws.Range("CA18").Value = 5738
$ws.Range("CB18").Value = "-"
$ws.Range("CC18").Value = "-"
$ws.Range("BY19").Value = 13280
$ws.Range("BZ19").Value = 6215
$ws.Range("CA19").Value = 7065
$ws.Range("CB19").Value = "-"
$ws.Range("CC19").Value = "-"
$ws.Range("BY20").Value = 16336
$ws.Range("BZ20").Value = 7655
$ws.Range("CA20").Value = 8681
$ws.Range("CB20").Value = "-"
$ws.Range("CC20").Value = "-"
$ws.Range("BY21").Value = 20529
$ws.Range("BZ21").Value = 9376
$ws.Range("CA21").Value = 11153
$ws.Range("CB21").Value = "-"
$ws.Range("CC21").Value = "-"
$ws.Range("BY22").Value = 25362
$ws.Range("BZ22").Value = 11658
$ws.Range("CA22").Value = 13704
$ws.Range("CB22").Value = "-"
$ws.Range("CC22").Value = "-"
$ws.Range("BY23").Value = 32582
$ws.Range("BZ23").Value = 15711
$ws.Range("CA23").Value = 16871
$ws.Range("CB23").Value = "-"
$ws.Range("CC23").Value = "-"
$ws.Range("BY24").Value = 39163
$ws.Range("BZ24").Value = 18767
$ws.Range("CA24").Value = 20396
$ws.Range("CB24").Value = "-"
$ws.Range("CC24").Value = "-"
$ws.Range("BY25").Value = 45066
$ws.Range("BZ25").Value = 21979
$ws.Range("CA25").Value = 23087
$ws.Range("CB25").Value = "-"
$ws.Range("CC25").Value = "-"
$ws.Range("BY26").Value = 49675
$ws.Range("BZ26").Value = 24293
$ws.Range("CA26").Value = 25382
$ws.Range("CB26").Value = "-"
$ws.Range("CC26").Value = "-"
$ws.Range("BY27").Value = 52115
$ws.Range("BZ27").Value = 25473
$ws.Range("CA27").Value = 26642
$ws.Range("CB27").Value = "-"
$ws.Range("CC27").Value = "-"
$ws.Range("BY28").Value = 55963
$ws.Range("BZ28").Value = 27579
$ws.Range("CA28").Value = 28384
$ws.Range("CB28").Value = "-"
$ws.Range("CC28").Value = "-"
$ws.Range("BY29").Value = 57966
$ws.Range("BZ29").Value = 29185
$ws.Range("CA29").Value = 28781
$ws.Range("CB29").Value = "-"
$ws.Range("CC29").Value = "-"
$ws.Range("CB30").Value = "-"
$ws.Range("CC30").Value = "-"
$ws.Range("BJ31").Value = 123734
$ws.Range("BK31").Value = 68139
$ws.Range("BL31").Value = 55595
$ws.Range("BM31").Value = "-"
$ws.Range("BN31").Value = "-"
$ws.Range("BO31").Value = 15095
$ws.Range("BP31").Value = 9185
$ws.Range("BQ31").Value = 5910
$ws.Range("BY31").Value = 13738
$ws.Range("BZ31").Value = 7701
$ws.Range("CA31").Value = 6037
$ws.Range("CB31").Value = "-"
$ws.Range("CC31").Value = "-"
$ws.Range("BJ32").Value = 132041
$ws.Range("BK32").Value = 73737
$ws.Range("BL32").Value = 58304
$ws.Range("BM32").Value = "-"
$ws.Range("BN32").Value = "-"
$ws.Range("BO32").Value = 16072
$ws.Range("BP32").Value = 9626
$ws.Range("BQ32").Value = 6446
$ws.Range("BR32").Value = "-"
$ws.Range("BS32").Value = "-"
$ws.Range("BY32").Value = 1179
$ws.Range("BZ32").Value = 606
$ws.Range("CA32").Value = 573
$ws.Range("CB32").Value = "-"
$ws.Range("CC32").Value = "-"
$ws.Range("BJ33").Value = 124141
$ws.Range("BK33").Value = 69279
$ws.Range("BL33").Value = 54862
$ws.Range("BM33").Value = "-"
$ws.Range("BN33").Value = "-"
$ws.Range("BO33").Value = 16894
$ws.Range("BP33").Value = 10270
$ws.Range("BQ33").Value = 6624
$ws.Range("BR33").Value = "-"
$ws.Range("BS33").Value = "-"
$ws.Range("BY33").Value = "-"
$ws.Range("BZ33").Value = "-"
$ws.Range("CA33").Value = "-"
$ws.Range("CB33").Value = "-"
$ws.Range("CC33").Value = "-"
$ws.Range("BJ34").Value = 117828
$ws.Range("BK34").Value = 65663
$ws.Range("BL34").Value = 52165
$ws.Range("BM34").Value = "-"
$ws.Range("BN34").Value = "-"
$ws.Range("BO34").Value = 15799
$ws.Range("BP34").Value = 9425
$ws.Range("BQ34").Value = 6374
$ws.Range("BR34").Value = "-"
$ws.Range("BS34").Value = "-"
$ws.Range("BY34").Value = "-"
$ws.Range("BZ34").Value = "-"
$ws.Range("CA34").Value = "-"
$ws.Range("CB34").Value = "-"
$ws.Range("CC34").Value = "-"
$ws.Range("BJ35").Value = 101067
$ws.Range("BK35").Value = 56587
$ws.Range("BL35").Value = 44480
$ws.Range("BM35").Value = "-"
$ws.Range("BN35").Value = "-"
$ws.Range("BO35").Value = 16473
$ws.Range("BP35").Value = 9985
$ws.Range("BQ35").Value = 6488
$ws.Range("BR35").Value = "-"
$ws.Range("BS35").Value = "-"
$ws.Range("BY35").Value = "-"
$ws.Range("BZ35").Value = "-"
$ws.Range("CA35").Value = "-"
$ws.Range("CB35").Value = "-"
$ws.Range("CC35").Value = "-"
$ws.Range("BJ36").Value = 97221
$ws.Range("BK36").Value = 54899
$ws.Range("BL36").Value = 42318
$ws.Range("BM36").Value = 4
$ws.Range("BN36").Value = "-"
$ws.Range("BO36").Value = 17157
$ws.Range("BP36").Value = 10482
$ws.Range("BQ36").Value = 6674
$ws.Range("BR36").Value = 1
$ws.Range("BS36").Value = "-"
$ws.Range("BY36").Value = "-"
$ws.Range("BZ36").Value = "-"
$ws.Range("CA36").Value = "-"
$ws.Range("CB36").Value = "-"
$ws.Range("CC36").Value = "-"

